$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 37 (pushes the Thyroid/Parathyroid rows down to 38/39) ---
$ws.Rows.Item(37).Insert() | Out-Null

# --- Populate the new row 37 with the "Splenic Artery Aneurysm" entry ---
$ws.Range("A37").Value = "Spleen"
$ws.Range("B37").Value = "Splenic Artery Aneurysm "
$ws.Range("C37").Value = "Clip 1 B-mode + Color Doppler"
$ws.Range("D37").Value = "https://youtu.be/nKZNmbVAutI"

# --- Rebuild the hyperlinks collection: same targets/order as before, but with ---
# --- the two refs that shifted down one row (old D37->D38, old D38->D39) fixed, ---
# --- then the brand-new hyperlink for row 37 appended last. ---
$links = @(
  @{ref="D4";  url="https://youtu.be/zxTC0YBY2RY"},
  @{ref="D29"; url="https://youtu.be/xBfd04F4Ni8"},
  @{ref="D12"; url="https://youtu.be/91M82AIMyu0"},
  @{ref="D35"; url="https://youtu.be/qushjTAy6XQ"},
  @{ref="D31"; url="https://youtu.be/pc-vbxSRTbs"},
  @{ref="D22"; url="https://youtu.be/DjI1kEnzfSQ"},
  @{ref="D30"; url="https://youtu.be/JvwODCASLYQ"},
  @{ref="D24"; url="https://youtu.be/U3ydTsRwxok"},
  @{ref="D15"; url="https://youtu.be/15o_Km86IzM"},
  @{ref="D36"; url="https://youtu.be/_FckFwJwynI"},
  @{ref="D32"; url="https://youtu.be/Axbee4vjNtU"},
  @{ref="D17"; url="https://youtu.be/RhSUFLTmTl4"},
  @{ref="D8";  url="https://youtu.be/2kRZcpi70Aw"},
  @{ref="D38"; url="https://youtu.be/z_oaRVxRz5s"},
  @{ref="D5";  url="https://youtu.be/K2Wbg7BgXy4"},
  @{ref="D3";  url="https://youtu.be/ZXwd0gwHEkQ"},
  @{ref="D33"; url="https://youtu.be/VJdnjrAAO-4"},
  @{ref="D2";  url="https://youtu.be/kdZO1IPuOIw"},
  @{ref="D39"; url="https://youtu.be/S45odD2wQOQ"},
  @{ref="D27"; url="https://youtu.be/ytNgK7wuL_M"},
  @{ref="D7";  url="https://youtu.be/mnDuOgdSpLA"},
  @{ref="D34"; url="https://youtu.be/aoaF345dsKc"},
  @{ref="D37"; url="https://youtu.be/nKZNmbVAutI"}
)

$ws.Hyperlinks.Delete() | Out-Null
foreach ($l in $links) {
    $ws.Hyperlinks.Add($ws.Range($l.ref), $l.url) | Out-Null
}

# --- Restore the D column hyperlink-style formatting that Hyperlinks.Delete() strips ---
$ws.Range("D2:D39").Style = "Collegamento ipertestuale"

# --- Reflect the selection state captured at save time ---
$null = $ws.Range("D45").Select()
